$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 55 (shifts existing rows 55..169 down to 56..170)
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with the new data record
$ws.Range("A55").Value = 8
$ws.Range("B55").Value = "Terminal La Palmera de La Serena"
$ws.Range("C55").Value = "Coquimbo"
$ws.Range("D55").Value = 45002
$ws.Range("E55").Value = 4
$ws.Range("F55").Value = 100112052
$ws.Range("G55").Value = "Albahaca"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 1040
$ws.Range("K55").Value = 2500
$ws.Range("L55").Value = 3000
$ws.Range("M55").Value = 2750
$ws.Range("N55").Value = "$/docena de matas"
$ws.Range("O55").Value = "Provincia del Elquí"
$ws.Range("P55").Value = 458
$ws.Range("Q55").Value = 6
$ws.Range("R55").Value = "Hortaliza"

# Match the date cell number format used by the other rows in column D
$ws.Range("D55").NumberFormat = $ws.Range("D56").NumberFormat
